$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: mark as Invalid (G) and Absent (H)
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Rows 4-6: mark Total Attendance Count (D) and Real (E)
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1

# Rows 7-18: mark as Absent (H)
for ($r = 7; $r -le 18; $r++) {
    $ws.Cells.Item($r, 8).Value = 1
}
